# Rename the worksheet from "EJERCICIO" to "2B"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "2B"
